$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '68.694.00'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '3.846.76'
$ws.Range("E3").Value = '  +2.96%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'600.88"
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").Value = "'163.86"
$ws.Range("E6").Value = '  -2.01%  '
$ws.Range("D7").Value = '3.844.93'
$ws.Range("E7").Value = '  +2.96%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -1.89%  '
$ws.Range("E10").Value = '  -1.17%  '
$ws.Range("D11").Value = "'6.35"
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").Value = "'36.98"
$ws.Range("E13").Value = '  -2.83%  '
$ws.Range("D14").Value = "'0.0000245"
$ws.Range("E14").Value = '  -1.05%  '
$ws.Range("D15").Value = '4.490.74'
$ws.Range("E15").Value = '  +2.97%  '
$ws.Range("D16").Value = '3.869.98'
$ws.Range("E16").Value = '  +3.43%  '
$ws.Range("D17").Value = '68.834.11'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("E18").Value = '  +2.85%  '
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("D20").Value = "'17.11"
$ws.Range("E20").Value = '  -1.24%  '
$ws.Range("D21").Value = "'11.26"
$ws.Range("E21").Value = '  +0.32%  '
$ws.Range("D22").Value = "'485.68"
$ws.Range("E22").Value = '  -1.18%  '
$ws.Range("E23").Value = '  -1.09%  '
$ws.Range("E24").Value = '  +6.84%  '
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = '  -1.59%  '
$ws.Range("D27").Value = "'12.13"
$ws.Range("E27").Value = '  -1.15%  '
$ws.Range("D28").Value = "'10.01"
$ws.Range("E28").Value = '  -0.60%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").Value = "'2.96"
$ws.Range("E30").Value = '  -0.59%  '
$ws.Range("E31").Value = '  -3.78%  '
$ws.Range("D32").Value = '3.997.32'
$ws.Range("E32").Value = '  +2.94%  '
$ws.Range("E33").Value = '  -3.27%  '
$ws.Range("E34").Value = '  +1.84%  '
$ws.Range("D35").Value = '3.792.98'
$ws.Range("E35").Value = '  +3.37%  '
$ws.Range("E36").Value = '  -0.98%  '
$ws.Range("E37").Value = '  +1.52%  '
$ws.Range("E38").Value = '  +0.24%  '
$ws.Range("D39").Value = "'5.89"
$ws.Range("E39").Value = '  -0.66%  '
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("E41").Value = '  -1.90%  '
$ws.Range("D42").Value = "'2.99"
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("D43").Value = "'433.19"
$ws.Range("E43").Value = '  +2.36%  '
$ws.Range("E44").Value = '  -0.55%  '
$ws.Range("D45").Value = "'1.99"
$ws.Range("E45").Value = '  +0.23%  '
$ws.Range("D47").Value = "'8.41"
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("D48").Value = "'142.98"
$ws.Range("E48").Value = '  +1.12%  '
$ws.Range("D49").Value = '2.841.79'
$ws.Range("E49").Value = '  +2.28%  '
$ws.Range("D50").Value = "'0.0358"
$ws.Range("E50").Value = '  +1.18%  '
$ws.Range("D51").Value = "'25.62"
$ws.Range("E51").Value = '  +12.77%  '
